# Append a new price-history row (row 17) to the sheet, mirroring the
# existing Datum/Cena rows, and keep the chart's source ranges in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 17

# A17: timestamp stored as text (matches existing rows, which are plain
# strings rather than real dates).
$ws.Cells.Item($newRow, 1).Value = "2026-01-24T20:37:37"

# B17: price value (numeric).
$ws.Cells.Item($newRow, 2).Value = 6636

# Extend the line chart's category/value series to include the new row.
$chartObject = $ws.ChartObjects().Item(1)
$chart = $chartObject.Chart
$series = $chart.SeriesCollection().Item(1)
$sheetName = $ws.Name
$series.XValues = "='" + $sheetName + "'!`$A`$9:`$A`$17"
$series.Values = "='" + $sheetName + "'!`$B`$9:`$B`$17"
